# Correct error on Week3 workbook first worksheet ("Logical Tests"):
# the logical-test descriptions in column C incorrectly referred to
# cells in column B, when the actual data being tested lives in
# column A. Update the description text so it references column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logical Tests")

$ws.Range("C2").Value  = "A2 is greater than 75"
$ws.Range("C3").Value  = "A3 is less than 50"
$ws.Range("C4").Value  = "A4 is equal to 100"
$ws.Range("C5").Value  = "A2 is greater than or equal to A3"
$ws.Range("C6").Value  = "A3 is less than or equal to A4"
$ws.Range("C7").Value  = "A4 is not equal to 50 (operator)"
$ws.Range("C8").Value  = "A8 is equal to ""apple"""
$ws.Range("C9").Value  = "A9 is equal to ""banana"""
$ws.Range("C10").Value = "A10 is equal to ""carrot"""
$ws.Range("C11").Value = "A11 is greater than 50 and less than 75"
$ws.Range("C12").Value = "A12 is greater than 50 or less than 75"
$ws.Range("C13").Value = "A2 is not equal to 99 (function)"
